# Add a new "18.0.0" release row (row 19) to the Versions sheet, with
# updated package version numbers for each project column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Give the new row's first cell (A19) the same bold style used by the
# other "Project / Release" cells in column A (copy format from A18).
$ws.Range("A18").Copy() | Out-Null
$ws.Range("A19").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Fill in the new values (order chosen to match how the shared-string
# table ends up being built, i.e. brand-new distinct strings first).
$ws.Range("A19").Value = "18.0.0"
$ws.Range("B19").Value = "15.1.0"
$ws.Range("I19").Value = "13.0.1"
$ws.Range("J19").Value = "3.0.2"
$ws.Range("C19").Value = "6.0.3"
$ws.Range("G19").Value = "10.1.1"
$ws.Range("M19").Value = "4.0.2"

# Remaining cells reuse version strings that already exist elsewhere in
# the workbook.
$ws.Range("D19").Value = "6.1.1"
$ws.Range("E19").Value = "6.0.1"
$ws.Range("F19").Value = "1.1.0"
$ws.Range("H19").Value = "13.0.0"
$ws.Range("K19").Value = "4.1.1"
$ws.Range("L19").Value = "6.1.1"

# Match the saved selection state from the edit (A19 selected).
$ws.Range("A19").Select() | Out-Null
